$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Seed the two new rows with the same formatting (date / integer styles)
# as the row directly above them, then overwrite the values.
$ws.Range("A69:F69").Copy()
$ws.Range("A70:F71").PasteSpecial(-4122)

# Row 70: 四方坪站 (2025-09-04)
$ws.Cells.Item(70, 1).Value = 45904
$ws.Cells.Item(70, 2).Value = "四方坪站"
$ws.Cells.Item(70, 3).Value = 11210.72
$ws.Cells.Item(70, 4).Value = 9183.56
$ws.Cells.Item(70, 5).Value = 3906.97
$ws.Cells.Item(70, 6).Value = 469

# Row 71: 高岭站 (2025-09-04)
$ws.Cells.Item(71, 1).Value = 45904
$ws.Cells.Item(71, 2).Value = "高岭站"
$ws.Cells.Item(71, 3).Value = 4480.51
$ws.Cells.Item(71, 4).Value = 3590.99
$ws.Cells.Item(71, 5).Value = 1085.17
$ws.Cells.Item(71, 6).Value = 159

$excel.CutCopyMode = $false

# Mirror the post-edit selected cell recorded in the worksheet view.
$ws.Range("G74").Select()
